# Remove the duplicate "pathology_distance_value" column (K) from the
# sample-block sheet and fix up the header-row cell comments so each
# comment once again documents the column it sits above (the comments
# do not automatically shift when the column is deleted, so we have to
# move their text manually), finally dropping the now orphaned comment
# for the column that fell off the end (old column T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original comment text for every header cell from L1 to T1
# (columns 12 through 20) before we touch anything - these texts need to
# slide one column to the left, into K1..S1.
$originalTexts = @{}
for ($col = 12; $col -le 20; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $originalTexts[$col] = $cell.Comment().Text()
}

# Delete the duplicate column K (pathology_distance_value). This shifts
# the cell values/data validations left by one column automatically, but
# leaves the comments behind attached to their old addresses.
$ws.Range("K1").EntireColumn().Delete()

# Re-apply the captured comment text, shifted one column to the left:
# what used to be the comment on column 12 (L) now belongs on column 11
# (K), column 13 (M) -> 12 (L), and so on through column 20 (T) -> 19 (S).
for ($col = 12; $col -le 20; $col++) {
    $targetCell = $ws.Cells.Item(1, $col - 1)
    [void]$targetCell.Comment().Text($originalTexts[$col])
}

# The last header cell (old column T / "Notes") no longer exists, so its
# comment (now duplicated at the previous step onto S1) must be removed.
[void]$ws.Cells.Item(1, 20).Comment().Delete()
